# DEV 5 - CLI Changes for Manager Project and Enquiry
#
# - Flip the "Melville Park" project (row 4) visibility from Hidden to Visible
# - Add a new BTO project "Archipelago" (Bedok Reservoir, row 6) as Visible

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Melville Park) visibility column (O) Hidden -> Visible
$ws.Cells.Item(4, 15).Value = "Visible"

# New row 6: Archipelago project at Bedok Reservoir
$ws.Cells.Item(6, 1).Value = 5          # Project ID
$ws.Cells.Item(6, 2).Value = "Archipelago"    # Project Name
$ws.Cells.Item(6, 3).Value = "Bedok Reservoir" # Neighborhood
$ws.Cells.Item(6, 4).Value = "2-ROOM"          # Type 1
$ws.Cells.Item(6, 5).Value = 10                # Number of units for Type 1
$ws.Cells.Item(6, 6).Value = 1000000           # Selling price for Type 1
$ws.Cells.Item(6, 7).Value = "3-ROOM"          # Type 2
$ws.Cells.Item(6, 8).Value = 10                # Number of units for Type 2
$ws.Cells.Item(6, 9).Value = 2000000           # Selling price for Type 2
$ws.Cells.Item(6, 10).Value = 45775            # Application opening date
$ws.Cells.Item(6, 11).Value = 45805            # Application closing date
$ws.Cells.Item(6, 12).Value = "T8765432F"      # Manager
$ws.Cells.Item(6, 13).Value = 10               # Officer Slot
$ws.Cells.Item(6, 15).Value = "Visible"        # Visible
